# Trade #71 closed at 2026-02-17 21:12:30 - unknown UNKNOWN +0.000%
#
# This script applies the new-trade / closed-trade update described by the
# commit. It:
#   1) Updates the Summary sheet roll-up metrics.
#   2) Updates the Strategy Status row for MarketMaking.
#   3) Updates the closed trade (row for Trade #99) on "All Trades" and
#      "MarketMaking" sheets, and appends a brand-new OPEN trade (Trade #132)
#      as a new row at the bottom of both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a date/time-looking string as real text (prevents Excel's
# COM layer from auto-coercing "2026-02-17" into a date serial number).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    # Reset to the default "Normal" style so the cell doesn't keep a
    # lingering Text number-format / quote-prefix style index once the
    # literal string value has already been committed.
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.11   # Current Capital
$summary.Range("B4").Value = 0.9       # Total P&L $
$summary.Range("B6").Value = 99        # Total Trades
$summary.Range("B7").Value = 47        # Winning Trades
$summary.Range("B9").Value = 47.47     # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.11
$status.Range("D5").Value = 66
$status.Range("E5").Value = 0.79
$status.Range("F5").Value = 1.11
$status.Range("G5").Value = 50

# ---------------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# 3a) Close out existing Trade #99 (sheet row 100)
$allTrades.Cells.Item(100, 7).Value = 0.86                 # Exit Price
$allTrades.Cells.Item(100, 8).Value = "CLOSED"              # Status
$allTrades.Cells.Item(100, 9).Value = 3.6145                # P&L %
$allTrades.Cells.Item(100, 10).Value = 0.03                 # P&L $
$allTrades.Cells.Item(100, 11).Value = 101.11                # Capital After
$allTrades.Cells.Item(100, 12).Value = "early_exit"          # Exit Reason
$allTrades.Cells.Item(100, 13).Value = 0.13                  # Duration (min)

# 3b) Append new Trade #132 (sheet row 133, brand new OPEN trade)
$allTrades.Cells.Item(133, 1).Value = 132
Set-TextValue $allTrades.Cells.Item(133, 2) "2026-02-17"
Set-TextValue $allTrades.Cells.Item(133, 3) "21:12:24"
$allTrades.Cells.Item(133, 4).Value = "MarketMaking"
$allTrades.Cells.Item(133, 5).Value = "UP"
$allTrades.Cells.Item(133, 6).Value = 0.83
$allTrades.Cells.Item(133, 7).Value = ""
$allTrades.Cells.Item(133, 8).Value = "OPEN"
$allTrades.Cells.Item(133, 9).Value = 0
$allTrades.Cells.Item(133, 10).Value = 0
$allTrades.Cells.Item(133, 11).Value = 101.0796151053151
$allTrades.Cells.Item(133, 12).Value = ""
$allTrades.Cells.Item(133, 13).Value = 0
$allTrades.Cells.Item(133, 14).Value = 0
$allTrades.Cells.Item(133, 15).Value = 0
$allTrades.Cells.Item(133, 16).Value = 0.6
$allTrades.Cells.Item(133, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# 4) MarketMaking sheet (strategy-specific trade log - same trades, different
#    column layout: L/M = Entry/Exit slippage, N = Confidence, O = Entry
#    Reason, P = Exit Reason, Q = Duration (min))
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# 4a) Close out existing Trade #99 (sheet row 67)
$mm.Cells.Item(67, 7).Value = 0.86                 # Exit Price
$mm.Cells.Item(67, 8).Value = "CLOSED"              # Status
$mm.Cells.Item(67, 9).Value = 3.6145                # P&L %
$mm.Cells.Item(67, 10).Value = 0.03                 # P&L $
$mm.Cells.Item(67, 11).Value = 101.11                # Capital After
$mm.Cells.Item(67, 16).Value = "early_exit"          # Exit Reason
$mm.Cells.Item(67, 17).Value = 0.13                  # Duration (min)

# 4b) Append new Trade #132 (sheet row 100, brand new OPEN trade)
$mm.Cells.Item(100, 1).Value = 132
Set-TextValue $mm.Cells.Item(100, 2) "2026-02-17"
Set-TextValue $mm.Cells.Item(100, 3) "21:12:24"
$mm.Cells.Item(100, 4).Value = "MarketMaking"
$mm.Cells.Item(100, 5).Value = "UP"
$mm.Cells.Item(100, 6).Value = 0.83
$mm.Cells.Item(100, 7).Value = ""
$mm.Cells.Item(100, 8).Value = "OPEN"
$mm.Cells.Item(100, 9).Value = 0
$mm.Cells.Item(100, 10).Value = 0
$mm.Cells.Item(100, 11).Value = 101.0796151053151
$mm.Cells.Item(100, 12).Value = 0
$mm.Cells.Item(100, 13).Value = 0
$mm.Cells.Item(100, 14).Value = 0.6
$mm.Cells.Item(100, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(100, 16).Value = ""
$mm.Cells.Item(100, 17).Value = 0
